$wb = $excel.ActiveWorkbook

# --- Sheet: trait_pixel ---
$ws1 = $wb.Worksheets.Item("trait_pixel")

# Header renames
$ws1.Range("D1").Value = "tassel_area ratio"
$ws1.Range("E1").Value = "average_width"
$ws1.Range("F1").Value = "average_height"
$ws1.Range("G1").Value = "number_branches"
$ws1.Range("H1").Value = "average_branch_length"
$ws1.Range("I1").Value = "average_coins_width"
$ws1.Range("K1").Value = "pixel/cm_ratio"

# Data updates (row 2)
$ws1.Range("G2").Value = 11
$ws1.Range("H2").Value = 587
$ws1.Range("I2").Value = 118
$ws1.Range("K2").Value = 43.7037037037037

# --- Sheet: trait_cm ---
$ws2 = $wb.Worksheets.Item("trait_cm")

# Header renames
$ws2.Range("E1").Value = "average_width"
$ws2.Range("F1").Value = "average_height"
$ws2.Range("G1").Value = "number_of_branches"
$ws2.Range("H1").Value = "average_branch_length"
$ws2.Range("I1").Value = "average_coins_width"
$ws2.Range("K1").Value = "pixel/cm_ratio"

# Data updates (row 2)
$ws2.Range("C2").Value = 125.1653314421143
$ws2.Range("E2").Value = 30.43220338983051
$ws2.Range("F2").Value = 46.60932203389831
$ws2.Range("G2").Value = 11
$ws2.Range("H2").Value = 13.43135593220339
$ws2.Range("K2").Value = 43.7037037037037
